# Applies the "GIT UPDATE" edit to the Rules worksheet:
#  - Cell E8 changes from "Good Morning" to "GIT UPDATE"
#  - Cell E8 becomes the active selection on the sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("E8").Value = "GIT UPDATE"

# Make E8 the active cell / selection, matching the recorded sheetView change.
$ws.Activate()
$ws.Range("E8").Select()
